$wb = $excel.ActiveWorkbook

# Rename the sheets (task order names embed new timestamps)
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556179750173"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556201205363"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255620126475"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556201824403"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651255620255486"

# Sheet 1 (GNG) - update stim file names
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512556179460182.csv"
$ws1.Range("B3").Value = "GNG_stims-16512556179580152.csv"
$ws1.Range("B4").Value = "go_stims-1651255617960025.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556179740167.csv"

# Sheet 2 (NB) - update stim file names
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512556186407766.csv"
$ws2.Range("B3").Value = "ZB-match_1-16512556181567786.csv"
$ws2.Range("B4").Value = "TB-16512556193897786.csv"
$ws2.Range("B5").Value = "TB-1651255620100569.csv"
$ws2.Range("B6").Value = "ZB-match_1-165125561810302.csv"
$ws2.Range("B7").Value = "TB-16512556194307778.csv"
$ws2.Range("B8").Value = "ZB-match_5-16512556180030143.csv"
$ws2.Range("B9").Value = "OB-16512556193707767.csv"
$ws2.Range("B10").Value = "OB-16512556182467766.csv"

# Sheet 4 (TOL) - update stim file names
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556201514308.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556201284804.csv"
$ws4.Range("B4").Value = "MM_stims-16512556201661777.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556201524327.csv"
$ws4.Range("B6").Value = "MM_stims-16512556201814427.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556201671808.csv"

# Sheet 5 (vSAT) - update stim file names
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512556201883333.csv"
$ws5.Range("B3").Value = "SAT_stims-16512556202079473.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512556202246907.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255620238901.csv"
